$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.485.06"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "1.921.69"
$ws.Range("E3").Value = "  +1.71%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "242.85"
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D8").Value = "0.2876"
$ws.Range("E8").Value = "  -3.32%  "
$ws.Range("D9").Value = "0.06781"
$ws.Range("E9").Value = "  +1.78%  "
$ws.Range("D10").Value = "107.27"
$ws.Range("E10").Value = "  +6.88%  "
$ws.Range("D11").Value = "18.32"
$ws.Range("E11").Value = "  -1.78%  "
$ws.Range("D12").Value = "0.07741"
$ws.Range("E12").Value = "  +2.40%  "
$ws.Range("D13").Value = "1.902.20"
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("D14").Value = "5.317"
$ws.Range("E14").Value = "  +3.17%  "
$ws.Range("D15").Value = "0.6585"
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("D16").Value = "293.75"
$ws.Range("E16").Value = "  -3.22%  "
$ws.Range("D17").Value = "30.482.67"
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("D18").Value = "0.000007615"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "12.97"
$ws.Range("E19").Value = "  -1.44%  "
$ws.Range("D20").Value = "0.9999"
$ws.Range("D21").Value = "2.146.08"
$ws.Range("E21").Value = "  +0.58%  "
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  +1.74%  "
$ws.Range("D24").Value = "6.194"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").Value = "9.376"
$ws.Range("E25").Value = "  +0.61%  "
$ws.Range("D26").Value = "168.85"
$ws.Range("E26").Value = "  +0.53%  "
$ws.Range("D27").Value = "21.50"
$ws.Range("E27").Value = "  +5.20%  "
$ws.Range("D28").Value = "2.125"
$ws.Range("E28").Value = "  +8.82%  "
$ws.Range("D29").Value = "0.1067"
$ws.Range("E29").Value = "  -5.51%  "
$ws.Range("D30").Value = "1.366"
$ws.Range("E30").Value = "  +1.00%  "
$ws.Range("D31").Value = "4.187"
$ws.Range("E31").Value = "  +0.27%  "
$ws.Range("D32").Value = "3.988"
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("D33").Value = "0.05035"
$ws.Range("E33").Value = "  -1.19%  "
$ws.Range("D34").Value = "0.7385"
$ws.Range("E34").Value = "  -1.97%  "
$ws.Range("D35").Value = "1.154"
$ws.Range("E35").Value = "  -0.82%  "
$ws.Range("D36").Value = "0.02096"
$ws.Range("E36").Value = "  +5.79%  "
$ws.Range("D37").Value = "2.739"
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("E38").Value = "  -0.98%  "
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("D40").Value = "110.32"
$ws.Range("E40").Value = "  +1.97%  "
$ws.Range("D41").Value = "0.8713"
$ws.Range("E41").Value = "  -2.82%  "
$ws.Range("D42").Value = "5.857"
$ws.Range("E42").Value = "  +3.67%  "
$ws.Range("D43").Value = "0.4253"
$ws.Range("E43").Value = "  +1.23%  "
$ws.Range("D44").Value = "1.0000"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").Value = "67.42"
$ws.Range("E45").Value = "  -0.41%  "
$ws.Range("D46").Value = "50.68"
$ws.Range("E46").Value = "  +18.52%  "
$ws.Range("D47").Value = "7.188"
$ws.Range("E47").Value = "  -2.62%  "
$ws.Range("D48").Value = "9.300"
$ws.Range("E48").Value = "  +2.52%  "
$ws.Range("D49").Value = "0.1216"
$ws.Range("E49").Value = "  -1.31%  "
$ws.Range("D50").Value = "35.01"
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").Value = "0.2465"
$ws.Range("E51").Value = "  +10.59%  "
